# "lab background, start and login art"
# - Renames the game title, adds a "begin"/BEGIN row under it.
# - Adds a "confirm"/CONFIRM row under "check".
# - Adds a 4-row "enter name" login block under "help".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Title text: "Cybrary Game" -> "Cybrary Quest"
$ws.Range("B3").Value = "Cybrary Quest"

# 2) Insert a new row right after "title" (row 3) for begin/BEGIN.
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "begin"
$ws.Range("B4").Value = "BEGIN"

# 3) Insert a new row right after "check" (now row 11) for confirm/CONFIRM.
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "confirm"
$ws.Range("B12").Value = "CONFIRM"

# 4) Insert four new rows right after "help" (now row 15) for the
#    login / enter-name block.
$ws.Range("A16:A19").EntireRow.Insert()
$ws.Range("A16").Value = "enter_name_title"
$ws.Range("B16").Value = "Login"
$ws.Range("A17").Value = "enter_name_name"
$ws.Range("B17").Value = "Name:"
$ws.Range("A18").Value = "enter_name_initial"
$ws.Range("B18").Value = "Initials:"
$ws.Range("A19").Value = "enter_name_placeholder"
$ws.Range("B19").Value = "Enter Text…"

# Restore tab selection/scroll focus near the new login block.
$ws.Range("B19").Select()
